# Apply scheduled-runner updates to market price / profit columns (H:N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 354.13333  # H28 was 217.24138
$ws.Cells.Item(28, 9).Value = 376.0909  # I28 was 242.57895
$ws.Cells.Item(28, 10).Value = 293.75  # J28 was 169.1
$ws.Cells.Item(28, 11).Value = 376.0909  # K28 was 242.57895
$ws.Cells.Item(28, 12).Value = 293.75  # L28 was 169.1
$ws.Cells.Item(28, 13).Value = 108.9091  # M28 was 242.42105
$ws.Cells.Item(28, 14).Value = -1263.75  # N28 was -1139.1
$ws.Cells.Item(111, 8).Value = 753.1875  # H111 was 649.05
$ws.Cells.Item(111, 9).Value = 470.06668  # I111 was 432.8889
$ws.Cells.Item(111, 10).Value = 5000  # J111 was 2594.5
$ws.Cells.Item(111, 11).Value = 1410.20004  # K111 was 1298.6667
$ws.Cells.Item(111, 12).Value = 15000  # L111 was 7783.5
$ws.Cells.Item(111, 13).Value = 1656.79996  # M111 was 1768.3333
$ws.Cells.Item(111, 14).Value = -21134  # N111 was -13917.5
$ws.Cells.Item(112, 8).Value = 1497.1072  # H112 was 1504.7778
$ws.Cells.Item(112, 10).Value = 1567.3914  # J112 was 1580
$ws.Cells.Item(112, 12).Value = 4702.174199999999  # L112 was 4740
$ws.Cells.Item(112, 14).Value = -6918.174199999999  # N112 was -6956
$ws.Cells.Item(137, 8).Value = 876.4231  # H137 was 936.8421
$ws.Cells.Item(137, 9).Value = 834.85  # I137 was 881.25
$ws.Cells.Item(137, 10).Value = 1015  # J137 was 1233.3334
$ws.Cells.Item(137, 11).Value = 2504.55  # K137 was 2643.75
$ws.Cells.Item(137, 12).Value = 3045  # L137 was 3700.0002
$ws.Cells.Item(137, 13).Value = 45.44999999999982  # M137 was -93.75
$ws.Cells.Item(137, 14).Value = -8145  # N137 was -8800.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(3, 8).Value = 0  # H3 was 1005
$ws.Cells.Item(3, 9).Value = 0  # I3 was 1005
$ws.Cells.Item(3, 11).Value = 0  # K3 was 1005
$ws.Cells.Item(3, 13).ClearContents()  # M3 was -890
$ws.Cells.Item(15, 8).Value = 19800  # H15 was 5006.5
$ws.Cells.Item(15, 10).Value = 19800  # J15 was 5006.5
$ws.Cells.Item(15, 12).Value = 19800  # L15 was 5006.5
$ws.Cells.Item(15, 14).Value = -20500  # N15 was -5706.5
$ws.Cells.Item(32, 8).Value = 6761.067  # H32 was 10420.381
$ws.Cells.Item(32, 9).Value = 6164.83  # I32 was 5955.5635
$ws.Cells.Item(32, 10).Value = 11275.429  # J32 was 41116
$ws.Cells.Item(32, 11).Value = 6164.83  # K32 was 5955.5635
$ws.Cells.Item(32, 12).Value = 11275.429  # L32 was 41116
$ws.Cells.Item(32, 13).Value = -5877.83  # M32 was -5668.5635
$ws.Cells.Item(32, 14).Value = -11849.429  # N32 was -41690
$ws.Cells.Item(61, 8).Value = 2500  # H61 was 1229.2222
$ws.Cells.Item(61, 9).Value = 0  # I61 was 1132.875
$ws.Cells.Item(61, 10).Value = 2500  # J61 was 2000
$ws.Cells.Item(61, 11).Value = 0  # K61 was 1132.875
$ws.Cells.Item(61, 12).Value = 2500  # L61 was 2000
$ws.Cells.Item(61, 13).ClearContents()  # M61 was -920.875
$ws.Cells.Item(61, 14).Value = -2924  # N61 was -2424
$ws.Cells.Item(63, 8).Value = 1944.75  # H63 was 1959.2
$ws.Cells.Item(63, 9).Value = 1944.75  # I63 was 1972.5
$ws.Cells.Item(63, 10).Value = 0  # J63 was 1906
$ws.Cells.Item(63, 11).Value = 1944.75  # K63 was 1972.5
$ws.Cells.Item(63, 12).Value = 0  # L63 was 1906
$ws.Cells.Item(63, 13).Value = -1258.75  # M63 was -1286.5
$ws.Cells.Item(63, 14).ClearContents()  # N63 was -3278
$ws.Cells.Item(66, 8).Value = 1944.75  # H66 was 1959.2
$ws.Cells.Item(66, 9).Value = 1944.75  # I66 was 1972.5
$ws.Cells.Item(66, 10).Value = 0  # J66 was 1906
$ws.Cells.Item(66, 11).Value = 9723.75  # K66 was 9862.5
$ws.Cells.Item(66, 12).Value = 0  # L66 was 9530
$ws.Cells.Item(66, 13).Value = -6291.75  # M66 was -6430.5
$ws.Cells.Item(66, 14).ClearContents()  # N66 was -16394
$ws.Cells.Item(74, 8).Value = 1065.9429  # H74 was 1230.6072
$ws.Cells.Item(74, 9).Value = 1080.5385  # I74 was 1328.579
$ws.Cells.Item(74, 11).Value = 1080.5385  # K74 was 1328.579
$ws.Cells.Item(74, 13).Value = -206.5385000000001  # M74 was -454.579
$ws.Cells.Item(77, 8).Value = 1065.9429  # H77 was 1230.6072
$ws.Cells.Item(77, 9).Value = 1080.5385  # I77 was 1328.579
$ws.Cells.Item(77, 11).Value = 5402.692500000001  # K77 was 6642.895
$ws.Cells.Item(77, 13).Value = -1034.692500000001  # M77 was -2274.895
$ws.Cells.Item(132, 8).Value = 3559.4  # H132 was 1682.3334
$ws.Cells.Item(132, 9).Value = 1400  # I132 was 987.3158
$ws.Cells.Item(132, 10).Value = 4999  # J132 was 2882.818
$ws.Cells.Item(132, 11).Value = 4200  # K132 was 2961.9474
$ws.Cells.Item(132, 12).Value = 14997  # L132 was 8648.454000000002
$ws.Cells.Item(132, 13).Value = -1670  # M132 was -431.9474
$ws.Cells.Item(132, 14).Value = -20057  # N132 was -13708.454
$ws.Cells.Item(136, 8).Value = 2500  # H136 was 1229.2222
$ws.Cells.Item(136, 9).Value = 0  # I136 was 1132.875
$ws.Cells.Item(136, 10).Value = 2500  # J136 was 2000
$ws.Cells.Item(136, 11).Value = 0  # K136 was 3398.625
$ws.Cells.Item(136, 12).Value = 7500  # L136 was 6000
$ws.Cells.Item(136, 13).ClearContents()  # M136 was -848.625
$ws.Cells.Item(136, 14).Value = -12600  # N136 was -11100

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 6356.4287  # H105 was 4822.6313
$ws.Cells.Item(105, 9).Value = 5415.8335  # I105 was 4803.3335
$ws.Cells.Item(105, 10).Value = 12000  # J105 was 4895
$ws.Cells.Item(105, 11).Value = 5415.8335  # K105 was 4803.3335
$ws.Cells.Item(105, 12).Value = 12000  # L105 was 4895
$ws.Cells.Item(105, 13).Value = -3668.8335  # M105 was -3056.3335
$ws.Cells.Item(105, 14).Value = -15494  # N105 was -8389
$ws.Cells.Item(134, 8).Value = 52672.78  # H134 was 30135.635
$ws.Cells.Item(134, 9).Value = 4510.05  # I134 was 3002.7925
$ws.Cells.Item(134, 10).Value = 98542.05  # J134 was 98613.766
$ws.Cells.Item(134, 11).Value = 13530.15  # K134 was 9008.377500000001
$ws.Cells.Item(134, 12).Value = 295626.15  # L134 was 295841.298
$ws.Cells.Item(134, 13).Value = -10995.15  # M134 was -6473.377500000001
$ws.Cells.Item(134, 14).Value = -300696.15  # N134 was -300911.298

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2738.963  # H31 was 2920.7693
$ws.Cells.Item(31, 9).Value = 2815.3044  # I31 was 3082.8572
$ws.Cells.Item(31, 10).Value = 2300  # J31 was 2240
$ws.Cells.Item(31, 11).Value = 2815.3044  # K31 was 3082.8572
$ws.Cells.Item(31, 12).Value = 2300  # L31 was 2240
$ws.Cells.Item(31, 13).Value = -2520.3044  # M31 was -2787.8572
$ws.Cells.Item(31, 14).Value = -2890  # N31 was -2830
$ws.Cells.Item(34, 8).Value = 2738.963  # H34 was 2920.7693
$ws.Cells.Item(34, 9).Value = 2815.3044  # I34 was 3082.8572
$ws.Cells.Item(34, 10).Value = 2300  # J34 was 2240
$ws.Cells.Item(34, 11).Value = 2815.3044  # K34 was 3082.8572
$ws.Cells.Item(34, 12).Value = 2300  # L34 was 2240
$ws.Cells.Item(34, 13).Value = -2613.3044  # M34 was -2880.8572
$ws.Cells.Item(34, 14).Value = -2704  # N34 was -2644
$ws.Cells.Item(58, 8).Value = 26978.5  # H58 was 8761.538
$ws.Cells.Item(58, 9).Value = 900  # I58 was 684.5
$ws.Cells.Item(58, 10).Value = 35671.332  # J58 was 21684.8
$ws.Cells.Item(58, 11).Value = 900  # K58 was 684.5
$ws.Cells.Item(58, 12).Value = 35671.332  # L58 was 21684.8
$ws.Cells.Item(58, 13).Value = -697  # M58 was -481.5
$ws.Cells.Item(58, 14).Value = -36077.332  # N58 was -22090.8
$ws.Cells.Item(86, 8).Value = 2400.5454  # H86 was 9777.346
$ws.Cells.Item(86, 9).Value = 2386.2727  # I86 was 18770.666
$ws.Cells.Item(86, 10).Value = 2414.818  # J86 was 2068.7856
$ws.Cells.Item(86, 11).Value = 2386.2727  # K86 was 18770.666
$ws.Cells.Item(86, 12).Value = 2414.818  # L86 was 2068.7856
$ws.Cells.Item(86, 13).Value = -1263.2727  # M86 was -17647.666
$ws.Cells.Item(86, 14).Value = -4660.818  # N86 was -4314.7856
$ws.Cells.Item(89, 8).Value = 2400.5454  # H89 was 9777.346
$ws.Cells.Item(89, 9).Value = 2386.2727  # I89 was 18770.666
$ws.Cells.Item(89, 10).Value = 2414.818  # J89 was 2068.7856
$ws.Cells.Item(89, 11).Value = 11931.3635  # K89 was 93853.33
$ws.Cells.Item(89, 12).Value = 12074.09  # L89 was 10343.928
$ws.Cells.Item(89, 13).Value = -6315.363499999999  # M89 was -88237.33
$ws.Cells.Item(89, 14).Value = -23306.09  # N89 was -21575.928
$ws.Cells.Item(132, 8).Value = 4084.6  # H132 was 2834.2354
$ws.Cells.Item(132, 9).Value = 3141.3333  # I132 was 2546.8667
$ws.Cells.Item(132, 10).Value = 5499.5  # J132 was 4989.5
$ws.Cells.Item(132, 11).Value = 9423.999899999999  # K132 was 7640.6001
$ws.Cells.Item(132, 12).Value = 16498.5  # L132 was 14968.5
$ws.Cells.Item(132, 13).Value = -6893.999899999999  # M132 was -5110.6001
$ws.Cells.Item(132, 14).Value = -21558.5  # N132 was -20028.5
$ws.Cells.Item(134, 8).Value = 6290.5557  # H134 was 4150.1875
$ws.Cells.Item(134, 9).Value = 7000.25  # I134 was 4032.625
$ws.Cells.Item(134, 10).Value = 5722.8  # J134 was 4267.75
$ws.Cells.Item(134, 11).Value = 21000.75  # K134 was 12097.875
$ws.Cells.Item(134, 12).Value = 17168.4  # L134 was 12803.25
$ws.Cells.Item(134, 13).Value = -18465.75  # M134 was -9562.875
$ws.Cells.Item(134, 14).Value = -22238.4  # N134 was -17873.25
$ws.Cells.Item(136, 8).Value = 26978.5  # H136 was 8761.538
$ws.Cells.Item(136, 9).Value = 900  # I136 was 684.5
$ws.Cells.Item(136, 10).Value = 35671.332  # J136 was 21684.8
$ws.Cells.Item(136, 11).Value = 2700  # K136 was 2053.5
$ws.Cells.Item(136, 12).Value = 107013.996  # L136 was 65054.39999999999
$ws.Cells.Item(136, 13).Value = -150  # M136 was 496.5
$ws.Cells.Item(136, 14).Value = -112113.996  # N136 was -70154.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 5347826.5  # H2 was 5347851
$ws.Cells.Item(2, 9).Value = 9804098  # I2 was 19608186
$ws.Cells.Item(2, 10).Value = 300.4  # J2 was 225.5
$ws.Cells.Item(2, 11).Value = 58824588  # K2 was 117649116
$ws.Cells.Item(2, 12).Value = 1802.4  # L2 was 1353
$ws.Cells.Item(2, 13).Value = -58824475  # M2 was -117649003
$ws.Cells.Item(2, 14).Value = -2028.4  # N2 was -1579

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(24, 8).Value = 8368337  # H24 was 7021002
$ws.Cells.Item(24, 10).Value = 7552505  # J24 was 70007
$ws.Cells.Item(24, 12).Value = 7552505  # L24 was 70007
$ws.Cells.Item(24, 14).Value = -7552851  # N24 was -70353
$ws.Cells.Item(70, 8).Value = 5386.75  # H70 was 5374.1875
$ws.Cells.Item(70, 9).Value = 4265.3335  # I70 was 4273.5
$ws.Cells.Item(70, 10).Value = 6828.5713  # J70 was 6474.875
$ws.Cells.Item(70, 11).Value = 4265.3335  # K70 was 4273.5
$ws.Cells.Item(70, 12).Value = 6828.5713  # L70 was 6474.875
$ws.Cells.Item(70, 13).Value = -3995.3335  # M70 was -4003.5
$ws.Cells.Item(70, 14).Value = -7368.5713  # N70 was -7014.875
$ws.Cells.Item(73, 8).Value = 5386.75  # H73 was 5374.1875
$ws.Cells.Item(73, 9).Value = 4265.3335  # I73 was 4273.5
$ws.Cells.Item(73, 10).Value = 6828.5713  # J73 was 6474.875
$ws.Cells.Item(73, 11).Value = 4265.3335  # K73 was 4273.5
$ws.Cells.Item(73, 12).Value = 6828.5713  # L73 was 6474.875
$ws.Cells.Item(73, 13).Value = -3329.3335  # M73 was -3337.5
$ws.Cells.Item(73, 14).Value = -8700.5713  # N73 was -8346.875
$ws.Cells.Item(132, 8).Value = 6666.3335  # H132 was 5727.75
$ws.Cells.Item(132, 9).Value = 11500  # I132 was 7206
$ws.Cells.Item(132, 11).Value = 34500  # K132 was 21618
$ws.Cells.Item(132, 13).Value = -31970  # M132 was -19088

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(43, 8).Value = 207000  # H43 was 253750
$ws.Cells.Item(43, 10).Value = 8750  # J43 was 5000
$ws.Cells.Item(43, 12).Value = 8750  # L43 was 5000
$ws.Cells.Item(43, 14).Value = -9136  # N43 was -5386
$ws.Cells.Item(132, 8).Value = 3068.121  # H132 was 2947.275
$ws.Cells.Item(132, 9).Value = 2603.0908  # I132 was 2705.1035
$ws.Cells.Item(132, 10).Value = 3998.182  # J132 was 3585.7273
$ws.Cells.Item(132, 11).Value = 7809.2724  # K132 was 8115.310500000001
$ws.Cells.Item(132, 12).Value = 11994.546  # L132 was 10757.1819
$ws.Cells.Item(132, 13).Value = -5279.2724  # M132 was -5585.310500000001
$ws.Cells.Item(132, 14).Value = -17054.546  # N132 was -15817.1819

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(3, 8).Value = 0  # H3 was 2352.5
$ws.Cells.Item(3, 9).Value = 0  # I3 was 2352.5
$ws.Cells.Item(3, 11).Value = 0  # K3 was 2352.5
$ws.Cells.Item(3, 13).ClearContents()  # M3 was -2238.5
$ws.Cells.Item(132, 8).Value = 2279  # H132 was 1647.3334
$ws.Cells.Item(132, 9).Value = 2045.5883  # I132 was 1398.68
$ws.Cells.Item(132, 10).Value = 3601.6667  # J132 was 2890.6
$ws.Cells.Item(132, 11).Value = 6136.7649  # K132 was 4196.04
$ws.Cells.Item(132, 12).Value = 10805.0001  # L132 was 8671.799999999999
$ws.Cells.Item(132, 13).Value = -3606.7649  # M132 was -1666.04
$ws.Cells.Item(132, 14).Value = -15865.0001  # N132 was -13731.8
$ws.Cells.Item(136, 8).Value = 1413.325  # H136 was 1400
$ws.Cells.Item(136, 9).Value = 1547.7858  # I136 was 1343.1082
$ws.Cells.Item(136, 10).Value = 1099.5834  # J136 was 1926.25
$ws.Cells.Item(136, 11).Value = 4643.357400000001  # K136 was 4029.3246
$ws.Cells.Item(136, 12).Value = 3298.7502  # L136 was 5778.75
$ws.Cells.Item(136, 13).Value = -2093.357400000001  # M136 was -1479.3246
$ws.Cells.Item(136, 14).Value = -8398.7502  # N136 was -10878.75

Write-Output "Applied $(211) cell updates across 8 sheets"